# Add placeholder statements to ddck
# - Update the "Automatic Connection Feature Completed" status text used throughout
#   column D to "Automatic Connection Feature Completed (Both inputs and outputs)"
# - Give row 37 (demands\dhw / WTap) its own, more specific status text:
#   "Automatic Connection Feature Completed, Updated port names"
# - Widen column D to fit the new, longer text
# - Move the active selection to D37

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldStatus = "Automatic Connection Feature Completed"
$newStatus = "Automatic Connection Feature Completed (Both inputs and outputs)"
$row37Status = "Automatic Connection Feature Completed, Updated port names"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# Row 37 gets its own distinct status message
$ws.Cells.Item(37, 4).Value = $row37Status

# Widen column D to fit the new text
$ws.Columns.Item(4).ColumnWidth = 56.5

# Update the active selection to D37
$ws.Range("D37").Select()
